$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($cellRef, $text)
    $ws.Range($cellRef).Value = "'" + $text
    $ws.Range($cellRef).Style = "Normal"
}

function Set-CellPlain {
    param($cellRef, $text)
    $ws.Range($cellRef).Value = $text
}

# Row 2
Set-CellPlain "D2" "69.471.77"
Set-CellPlain "E2" "  +1.93%  "

# Row 3
Set-CellPlain "D3" "3.910.39"
Set-CellPlain "E3" "  +0.53%  "

# Row 4
Set-CellPlain "E4" "  +0.15%  "

# Row 5
Set-CellText "D5" "530.68"
Set-CellPlain "E5" "  +9.87%  "

# Row 6
Set-CellText "D6" "144.91"
Set-CellPlain "E6" "  +0.16%  "

# Row 7
Set-CellPlain "E7" "  -1.31%  "

# Row 8
Set-CellPlain "E8" "  +0.04%  "

# Row 9
Set-CellText "D9" "0.719"
Set-CellPlain "E9" "  -2.76%  "

# Row 10
Set-CellPlain "E10" "  -5.86%  "

# Row 11
Set-CellPlain "E11" "  -6.13%  "

# Row 12
Set-CellText "D12" "42.19"
Set-CellPlain "E12" "  -1.58%  "

# Row 13
Set-CellPlain "D13" "4.531.33"
Set-CellPlain "E13" "  +0.62%  "

# Row 14
Set-CellPlain "E14" "  -2.73%  "

# Row 15
Set-CellPlain "D15" "3.927.08"
Set-CellPlain "E15" "  +0.67%  "

# Row 16
Set-CellText "D16" "14.06"
Set-CellPlain "E16" "  -1.43%  "

# Row 17
Set-CellPlain "E17" "  +8.77%  "

# Row 18
Set-CellPlain "E18" "  -0.64%  "

# Row 19
Set-CellText "D19" "19.81"
Set-CellPlain "E19" "  -0.92%  "

# Row 20
Set-CellPlain "D20" "69.413.72"
Set-CellPlain "E20" "  +1.83%  "

# Row 21
Set-CellText "D21" "429.84"
Set-CellPlain "E21" "  -0.04%  "

# Row 22
Set-CellPlain "E22" "  -5.08%  "

# Row 23
Set-CellText "D23" "14.21"
Set-CellPlain "E23" "  -4.35%  "

# Row 24
Set-CellText "D24" "88.55"
Set-CellPlain "E24" "  -1.39%  "

# Row 25
Set-CellPlain "E25" "  +9.48%  "

# Row 26
Set-CellText "D26" "11.52"
Set-CellPlain "E26" "  -3.96%  "

# Row 27
Set-CellText "D27" "10.63"
Set-CellPlain "E27" "  -3.49%  "

# Row 28
Set-CellText "D28" "36.44"
Set-CellPlain "E28" "  -2.65%  "

# Row 29
Set-CellText "D29" "688.93"
Set-CellPlain "E29" "  -3.09%  "

# Row 30
Set-CellText "D30" "13.18"
Set-CellPlain "E30" "  -2.73%  "

# Row 31
Set-CellPlain "E31" "  -3.12%  "

# Row 32
Set-CellText "D32" "2.83"
Set-CellPlain "E32" "  -2.34%  "

# Row 33
Set-CellText "D33" "68.27"
Set-CellPlain "E33" "  +11.98%  "

# Row 34
Set-CellText "D34" "0.444"
Set-CellPlain "E34" "  +11.64%  "

# Row 35
Set-CellPlain "E35" "  -1.57%  "

# Row 36
Set-CellText "D36" "40.08"
Set-CellPlain "E36" "  -2.17%  "

# Row 37
Set-CellPlain "E37" "  -3.17%  "

# Row 38
Set-CellPlain "E38" "  +2.68%  "

# Row 39
Set-CellText "D39" "1.00"
Set-CellPlain "E39" "  +0.05%  "

# Row 40
Set-CellPlain "E40" "  -0.14%  "

# Row 41
Set-CellPlain "E41" "  -4.29%  "

# Row 42
Set-CellText "D42" "3.17"
Set-CellPlain "E42" "  +2.85%  "

# Row 43
Set-CellPlain "E43" "  +5.86%  "

# Row 44
Set-CellText "D44" "2.80"
Set-CellPlain "E44" "  -5.63%  "

# Row 45
Set-CellPlain "B45" "ApeXProtocol"
Set-CellPlain "C45" "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-CellText "D45" "3.35"
Set-CellPlain "E45" "  -0.04%  "

# Row 46
Set-CellPlain "B46" "Stellar"
Set-CellPlain "C46" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-CellText "D46" "0.140"
Set-CellPlain "E46" "  -1.55%  "

# Row 47
Set-CellPlain "D47" "0.0₆0358"
Set-CellPlain "E47" "  +12.28%  "

# Row 48
Set-CellPlain "E48" "  +7.36%  "

# Row 49
Set-CellPlain "D49" "2.762.10"
Set-CellPlain "E49" "  +13.86%  "

# Row 50
Set-CellText "D50" "144.62"
Set-CellPlain "E50" "  +0.01%  "

# Row 51
Set-CellPlain "E51" "  -3.18%  "
